$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "Claudinei-Elemaq.-"
$ws.Range("C18").Value = "Claudinei-Des. Maq. Cad-"
$ws.Range("D18").Value = "Claudinei-Des. Maq. Cad-"
$ws.Range("E18").Value = "[Leonardo-Retífica-2NA, Leonardo-Retífica-2NA, Leonardo-Retífica-2NA, Leonardo-Retífica-2NA]"
$ws.Range("F18").Value = "[leonardo-M.Maq.E.I.-2NA, leonardo-M.Maq.E.I.-2NA, leonardo-M.Maq.E.I.-2NA, leonardo-M.Maq.E.I.-2NA]"

# Row 19
$ws.Range("B19").Value = "Claudinei-Elemaq.-"
$ws.Range("C19").Value = "[Euclides-Soldagem-2NA, Euclides-Soldagem-2NA, Euclides-Soldagem-2NA, Euclides-Soldagem-2NA]"
$ws.Range("D19").Value = "Claudinei-Des. Maq. Cad-"
$ws.Range("E19").Value = "[Paulo Rob.-Usin. CNC-2NA, Paulo Rob.-Usin. CNC-2NA, Paulo Rob.-Usin. CNC-2NA, Paulo Rob.-Usin. CNC-2NA]"
$ws.Range("F19").Value = "[João Bosco-Fundição-2NA, João Bosco-Fundição-2NA, João Bosco-Fundição-2NA, João Bosco-Fundição-2NA]"

# Row 20
$ws.Range("B20").Value = "Euclides-Gest. Int.-"
$ws.Range("C20").Value = "[Guilherme-C. Hidráulica-2NA, Guilherme-C. Hidráulica-2NA, Guilherme-C. Hidráulica-2NA, Guilherme-C. Hidráulica-2NA]"
$ws.Range("D20").Value = "[Ismail-Metrologia 2-2NA, Ismail-Metrologia 2-2NA, Ismail-Metrologia 2-2NA, Ismail-Metrologia 2-2NA]"
$ws.Range("E20").Value = "[Aderci-Fresagem-2NA, Aderci-Fresagem-2NA, Aderci-Fresagem-2NA, Aderci-Fresagem-2NA]"
$ws.Range("F20").Value = "[Guilherme-C.L.P.-2NA, Guilherme-C.L.P.-2NA, Guilherme-C.L.P.-2NA, Guilherme-C.L.P.-2NA]"

# Row 21
$ws.Range("B21").Value = "Euclides-Gest. Int.-"
$ws.Range("D21").Value = "[Leandro-M.S.R.A.C.-2NA, Leandro-M.S.R.A.C.-2NA, Leandro-M.S.R.A.C.-2NA, Leandro-M.S.R.A.C.-2NA]"
$ws.Range("E21").Value = "[Guilherme-C. Pneumática-2NA, Guilherme-C. Pneumática-2NA, Guilherme-C. Pneumática-2NA, Guilherme-C. Pneumática-2NA]"
$ws.Range("F21").Value = "[Aderci-CAD / CAM-2NA, Aderci-CAD / CAM-2NA, Aderci-CAD / CAM-2NA, Aderci-CAD / CAM-2NA]"
